# Apply updated crypto price/volume data to match source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.797.45"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "2.313.34"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.63%  "
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "2.664.10"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "2.313.18"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "42.748.49"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +32.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0890"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.50%  "
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0353"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "115.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "1.619.11"
$ws.Range("E51").Value = "  +5.39%  "
